$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Rspo2"
$ws.Cells.Item(2, 3).Value = "Znrf3"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.597878666666666
$ws.Cells.Item(2, 8).Value = 4.793635999999999
$ws.Cells.Item(2, 9).Value = 0.98224549682877
$ws.Cells.Item(2, 10).Value = 0.9822454968287699
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.8289893333333334
$ws.Cells.Item(2, 14).Value = 2.486968
$ws.Cells.Item(2, 15).Value = 0.06741872448909192
$ws.Cells.Item(2, 16).Value = 0.06741872448909192
$ws.Cells.Item(2, 17).Value = 1.324624370627555
$ws.Cells.Item(2, 18).Value = 11.921619335648
$ws.Cells.Item(2, 19).Value = 0.06622173853135006
$ws.Cells.Item(2, 20).Value = 0.06622173853135005

$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Rspo2"
$ws.Cells.Item(3, 3).Value = "Znrf3"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.597878666666666
$ws.Cells.Item(3, 8).Value = 4.793635999999999
$ws.Cells.Item(3, 9).Value = 0.98224549682877
$ws.Cells.Item(3, 10).Value = 0.9822454968287699
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 2.848096333333333
$ws.Cells.Item(3, 14).Value = 8.544289
$ws.Cells.Item(3, 15).Value = 0.2316254435305073
$ws.Cells.Item(3, 16).Value = 0.2316254435305073
$ws.Cells.Item(3, 17).Value = 4.550912371644888
$ws.Cells.Item(3, 18).Value = 40.95821134480399
$ws.Cells.Item(3, 19).Value = 0.2275130488588073
$ws.Cells.Item(3, 20).Value = 0.2275130488588073

$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Rspo2"
$ws.Cells.Item(4, 3).Value = "Znrf3"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.597878666666666
$ws.Cells.Item(4, 8).Value = 4.793635999999999
$ws.Cells.Item(4, 9).Value = 0.98224549682877
$ws.Cells.Item(4, 10).Value = 0.9822454968287699
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 1.875312666666667
$ws.Cells.Item(4, 14).Value = 5.625938
$ws.Cells.Item(4, 15).Value = 0.1525124424659717
$ws.Cells.Item(4, 16).Value = 0.1525124424659717
$ws.Cells.Item(4, 17).Value = 2.996522103396444
$ws.Cells.Item(4, 18).Value = 26.968698930568
$ws.Cells.Item(4, 19).Value = 0.1498046598225576
$ws.Cells.Item(4, 20).Value = 0.1498046598225576

$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Rspo2"
$ws.Cells.Item(5, 3).Value = "Znrf3"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.597878666666666
$ws.Cells.Item(5, 8).Value = 4.793635999999999
$ws.Cells.Item(5, 9).Value = 0.98224549682877
$ws.Cells.Item(5, 10).Value = 0.9822454968287699
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 6.743730666666667
$ws.Cells.Item(5, 14).Value = 20.231192
$ws.Cells.Item(5, 15).Value = 0.5484433895144291
$ws.Cells.Item(5, 16).Value = 0.5484433895144291
$ws.Cells.Item(5, 17).Value = 10.77566336601244
$ws.Cells.Item(5, 18).Value = 96.98097029411198
$ws.Cells.Item(5, 19).Value = 0.5387060496160551
$ws.Cells.Item(5, 20).Value = 0.538706049616055

$ws.Cells.Item(6, 1).Value = "sCs"
$ws.Cells.Item(6, 2).Value = "Rspo2"
$ws.Cells.Item(6, 3).Value = "Znrf3"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.02888233333333333
$ws.Cells.Item(6, 8).Value = 0.086647
$ws.Cells.Item(6, 9).Value = 0.01775450317123003
$ws.Cells.Item(6, 10).Value = 0.01775450317123003
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.8289893333333334
$ws.Cells.Item(6, 14).Value = 2.486968
$ws.Cells.Item(6, 15).Value = 0.06741872448909192
$ws.Cells.Item(6, 16).Value = 0.06741872448909192
$ws.Cells.Item(6, 17).Value = 0.02394314625511111
$ws.Cells.Item(6, 18).Value = 0.215488316296
$ws.Cells.Item(6, 19).Value = 0.001196985957741866
$ws.Cells.Item(6, 20).Value = 0.001196985957741866

$ws.Cells.Item(7, 1).Value = "sCs"
$ws.Cells.Item(7, 2).Value = "Rspo2"
$ws.Cells.Item(7, 3).Value = "Znrf3"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 0.3333333333333333
$ws.Cells.Item(7, 7).Value = 0.02888233333333333
$ws.Cells.Item(7, 8).Value = 0.086647
$ws.Cells.Item(7, 9).Value = 0.01775450317123003
$ws.Cells.Item(7, 10).Value = 0.01775450317123003
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 2.848096333333333
$ws.Cells.Item(7, 14).Value = 8.544289
$ws.Cells.Item(7, 15).Value = 0.2316254435305073
$ws.Cells.Item(7, 16).Value = 0.2316254435305073
$ws.Cells.Item(7, 17).Value = 0.08225966766477777
$ws.Cells.Item(7, 18).Value = 0.7403370089829999
$ws.Cells.Item(7, 19).Value = 0.004112394671699954
$ws.Cells.Item(7, 20).Value = 0.004112394671699954

$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Rspo2"
$ws.Cells.Item(8, 3).Value = "Znrf3"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 0.3333333333333333
$ws.Cells.Item(8, 7).Value = 0.02888233333333333
$ws.Cells.Item(8, 8).Value = 0.086647
$ws.Cells.Item(8, 9).Value = 0.01775450317123003
$ws.Cells.Item(8, 10).Value = 0.01775450317123003
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 1.875312666666667
$ws.Cells.Item(8, 14).Value = 5.625938
$ws.Cells.Item(8, 15).Value = 0.1525124424659717
$ws.Cells.Item(8, 16).Value = 0.1525124424659717
$ws.Cells.Item(8, 17).Value = 0.05416340554288888
$ws.Cells.Item(8, 18).Value = 0.487470649886
$ws.Cells.Item(8, 19).Value = 0.002707782643414133
$ws.Cells.Item(8, 20).Value = 0.002707782643414133

$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Rspo2"
$ws.Cells.Item(9, 3).Value = "Znrf3"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 0.3333333333333333
$ws.Cells.Item(9, 7).Value = 0.02888233333333333
$ws.Cells.Item(9, 8).Value = 0.086647
$ws.Cells.Item(9, 9).Value = 0.01775450317123003
$ws.Cells.Item(9, 10).Value = 0.01775450317123003
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 6.743730666666667
$ws.Cells.Item(9, 14).Value = 20.231192
$ws.Cells.Item(9, 15).Value = 0.5484433895144291
$ws.Cells.Item(9, 16).Value = 0.5484433895144291
$ws.Cells.Item(9, 17).Value = 0.1947746770248889
$ws.Cells.Item(9, 18).Value = 1.752972093224
$ws.Cells.Item(9, 19).Value = 0.009737339898374078
$ws.Cells.Item(9, 20).Value = 0.009737339898374078

